$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update shared string "very hard" -> "very easy" (AO column, M_node_feat_type entries etc. not affected)
$ws.Range("AO2:AO13").Value = "very easy"

# Update per-row toy data values
$ws.Range("A2").Value = 31.22534311612447
$ws.Range("B2").Value = 0.2828282828282828
$ws.Range("C2").Value = 3.89
$ws.Range("D2").Value = 17.38
$ws.Range("O2").Value = 76
$ws.Range("P2").Value = 305
$ws.Range("Q2").Value = 83
$ws.Range("R2").Value = 83

$ws.Range("A3").Value = 31.61318896214167
$ws.Range("B3").Value = 0.2222222222222222
$ws.Range("C3").Value = 3.91
$ws.Range("D3").Value = 16.56
$ws.Range("O3").Value = 76
$ws.Range("P3").Value = 305
$ws.Range("Q3").Value = 83
$ws.Range("R3").Value = 83

$ws.Range("A4").Value = 33.53106116453807
$ws.Range("B4").Value = 0.3131313131313131
$ws.Range("C4").Value = 3.68
$ws.Range("D4").Value = 15.72
$ws.Range("O4").Value = 76
$ws.Range("P4").Value = 305
$ws.Range("Q4").Value = 83
$ws.Range("R4").Value = 83

$ws.Range("A5").Value = 20.50241411526998
$ws.Range("B5").Value = 0.202020202020202
$ws.Range("C5").Value = 3.95
$ws.Range("D5").Value = 18.25
$ws.Range("O5").Value = 76
$ws.Range("P5").Value = 305
$ws.Range("Q5").Value = 83
$ws.Range("R5").Value = 83

$ws.Range("A6").Value = 20.97510018348694
$ws.Range("B6").Value = 0.2626262626262627
$ws.Range("C6").Value = 3.79
$ws.Range("D6").Value = 17.72
$ws.Range("O6").Value = 76
$ws.Range("P6").Value = 305
$ws.Range("Q6").Value = 83
$ws.Range("R6").Value = 83

$ws.Range("A7").Value = 21.01098755598068
$ws.Range("B7").Value = 0.2323232323232323
$ws.Range("C7").Value = 3.85
$ws.Range("D7").Value = 17.21
$ws.Range("O7").Value = 76
$ws.Range("P7").Value = 305
$ws.Range("Q7").Value = 83
$ws.Range("R7").Value = 83

$ws.Range("A8").Value = 11.10944559176763
$ws.Range("B8").Value = 0.2626262626262627
$ws.Range("C8").Value = 3.87
$ws.Range("D8").Value = 16.46
$ws.Range("O8").Value = 76
$ws.Range("P8").Value = 305
$ws.Range("Q8").Value = 83
$ws.Range("R8").Value = 83

$ws.Range("A9").Value = 11.95082873503367
$ws.Range("B9").Value = 0.3434343434343434
$ws.Range("C9").Value = 3.9
$ws.Range("D9").Value = 15.49
$ws.Range("O9").Value = 76
$ws.Range("P9").Value = 305
$ws.Range("Q9").Value = 83
$ws.Range("R9").Value = 83

$ws.Range("A10").Value = 11.89163700342178
$ws.Range("B10").Value = 0.4141414141414141
$ws.Range("C10").Value = 3.75
$ws.Range("D10").Value = 14.42
$ws.Range("O10").Value = 76
$ws.Range("P10").Value = 305
$ws.Range("Q10").Value = 83
$ws.Range("R10").Value = 83

$ws.Range("A11").Value = 8.372995054721832
$ws.Range("B11").Value = 0.2424242424242424
$ws.Range("C11").Value = 3.82
$ws.Range("D11").Value = 16.87
$ws.Range("O11").Value = 76
$ws.Range("P11").Value = 305
$ws.Range("Q11").Value = 83
$ws.Range("R11").Value = 83

$ws.Range("A12").Value = 8.535145397981008
$ws.Range("B12").Value = 0.1717171717171717
$ws.Range("C12").Value = 3.91
$ws.Range("D12").Value = 17.42
$ws.Range("O12").Value = 76
$ws.Range("P12").Value = 305
$ws.Range("Q12").Value = 83
$ws.Range("R12").Value = 83

$ws.Range("A13").Value = 8.086561504999796
$ws.Range("B13").Value = 0.2323232323232323
$ws.Range("C13").Value = 3.79
$ws.Range("D13").Value = 16.81
$ws.Range("O13").Value = 76
$ws.Range("P13").Value = 305
$ws.Range("Q13").Value = 83
$ws.Range("R13").Value = 83
